$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full "species observation" data between row pairs
# (3 <-> 5) and (4 <-> 6). Only the columns below actually differ between
# the two rows in each pair; every other column already holds an identical
# value in both rows of the pair, so leaving them untouched keeps the sheet
# correct while avoiding any unwanted side effects (e.g. Excel re-parsing
# untouched date-like text cells).
$cols = @("A","B","E","F","G","H","Q","R","AJ","AK","AO","AX")

function Swap-RowValues($rowA, $rowB, $columns) {
    foreach ($col in $columns) {
        $cellA = $ws.Range($col + $rowA)
        $cellB = $ws.Range($col + $rowB)

        $valueA = $cellA.Value2
        $valueB = $cellB.Value2

        if ($valueA -eq $null) { $valueA = "" }
        if ($valueB -eq $null) { $valueB = "" }

        $cellA.Value2 = $valueB
        $cellB.Value2 = $valueA
    }
}

Swap-RowValues 3 5 $cols
Swap-RowValues 4 6 $cols
